# Pressure NOR file changes
#
# 1) Re-order the "Output_*" sheets: the six "Manualadjust"/"Mnladjst_ntv"
#    sheets (KPA/BAR/PSI) move from the tail of the workbook to right after
#    the Input sheets, ahead of the electronicadjust/Runup output sheets.
# 2) Fix the data rows on a handful of output sheets whose values were
#    copied from the wrong source sheet (BAR/PSI Manualadjust + native
#    variants, and KPA electronicadjust + native variant).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: reorder sheets into the target tab order.
# ---------------------------------------------------------------------
$targetOrder = @(
    "Input_KPA_manualadjust",
    "Input_BAR_manualadjust",
    "Input_PSI_manualadjust",
    "Input_KPA_electronicadjust",
    "Input_BAR_electronicadjust",
    "Input_PSI_electronicadjust",
    "Input_KPA_Runup",
    "Input_BAR_Runup",
    "Input_PSI_Runup",
    "Input_Runup_LineSpeedft_min",
    "Input_Runup_LineSpeedm_min",
    "Output_KPA_Manualadjust",
    "Output_KPA_Mnladjst_ntv",
    "Output_BAR_Manualadjust",
    "Output_BAR_Mnladjst_ntv",
    "Output_PSI_Manualadjust",
    "Output_PSI_Mnladjst_ntv",
    "Output_KPA_electronicadjust",
    "Output_KPA_electronicadjust_Ntv",
    "Output_BAR_electronicadjust",
    "Output_BAR_electronicadjust_Ntv",
    "Output_PSI_electronicadjust",
    "Output_PSI_electronicadjust_Ntv",
    "Output_KPA_Runup",
    "Output_BAR_Runup",
    "Output_PSI_Runup",
    "Output_Runup_Lnspeed_ftmin",
    "Output_Runup_Lnspeed_mmin"
)

for ($i = 1; $i -lt $targetOrder.Length; $i++) {
    $prevSheet = $wb.Worksheets.Item($targetOrder[$i - 1])
    $ws = $wb.Worksheets.Item($targetOrder[$i])
    $ws.Move($null, $prevSheet)
}

# ---------------------------------------------------------------------
# Step 2: correct the cell values that were wrong in the source sheets.
# Values are stored as text (quote-prefixed) to match the original
# shared-string / text cell typing rather than becoming numbers.
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Output_BAR_Manualadjust")
$ws.Range("A2").Value = "'4080"
$ws.Range("B2").Value = "'82020"

$ws = $wb.Worksheets.Item("Output_BAR_Mnladjst_ntv")
$ws.Range("A2").Value = "'4931"
$ws.Range("B2").Value = "'99133"

$ws = $wb.Worksheets.Item("Output_PSI_Manualadjust")
$ws.Range("A2").Value = "'1140000"
$ws.Range("B2").Value = "'1200000"

$ws = $wb.Worksheets.Item("Output_PSI_Mnladjst_ntv")
$ws.Range("A2").Value = "'95000"
$ws.Range("B2").Value = "'100000"

$ws = $wb.Worksheets.Item("Output_KPA_electronicadjust")
$ws.Range("A2").Value = "'6796000"
$ws.Range("B2").Value = "'408000"
$ws.Range("C2").Value = "'4796000"
$ws.Range("D2").Value = "'5796000"
$ws.Range("E2").Value = "'8280000"

$ws = $wb.Worksheets.Item("Output_KPA_electronicadjust_Ntv")
$ws.Range("A2").Value = "'82140"
$ws.Range("B2").Value = "'4931"
$ws.Range("C2").Value = "'57967"
$ws.Range("D2").Value = "'70053"
$ws.Range("E2").Value = "'100076"
